{"js": "// Replace each old division-problem text with its new value.\n// Every cell text in the table is unique, so a direct search+replace\n// (matchCase + matchWholeWord) for each pair is unambiguous.\nconst replacements = [\n  [\n    \"473\u00f75=94, 3\",\n    \"162\u00f79=18, 0\"\n  ],\n  [\n    \"758\u00f72=379, 0\",\n    \"855\u00f75=171, 0\"\n  ],\n  [\n    \"310\u00f78=38, 6\",\n    \"365\u00f78=45, 5\"\n  ],\n  [\n    \"128\u00f77=18, 2\",\n    \"965\u00f73=321, 2\"\n  ],\n  [\n    \"607\u00f76=101, 1\",\n    \"421\u00f78=52, 5\"\n  ],\n  [\n    \"971\u00f79=107, 8\",\n    \"754\u00f79=83, 7\"\n  ],\n  [\n    \"698\u00f75=139, 3\",\n    \"175\u00f76=29, 1\"\n  ],\n  [\n    \"240\u00f73=80, 0\",\n    \"374\u00f79=41, 5\"\n  ],\n  [\n    \"764\u00f72=382, 0\",\n    \"572\u00f78=71, 4\"\n  ],\n  [\n    \"919\u00f79=102, 1\",\n    \"495\u00f79=55, 0\"\n  ],\n  [\n    \"219\u00f74=54, 3\",\n    \"567\u00f75=113, 2\"\n  ],\n  [\n    \"538\u00f77=76, 6\",\n    \"196\u00f78=24, 4\"\n  ],\n  [\n    \"399\u00f75=79, 4\",\n    \"744\u00f79=82, 6\"\n  ],\n  [\n    \"972\u00f75=194, 2\",\n    \"810\u00f79=90, 0\"\n  ],\n  [\n    \"770\u00f75=154, 0\",\n    \"341\u00f77=48, 5\"\n  ],\n  [\n    \"701\u00f76=116, 5\",\n    \"928\u00f76=154, 4\"\n  ],\n  [\n    \"421\u00f74=105, 1\",\n    \"838\u00f79=93, 1\"\n  ],\n  [\n    \"987\u00f78=123, 3\",\n    \"778\u00f73=259, 1\"\n  ],\n  [\n    \"237\u00f73=79, 0\",\n    \"682\u00f78=85, 2\"\n  ],\n  [\n    \"326\u00f76=54, 2\",\n    \"489\u00f78=61, 1\"\n  ],\n  [\n    \"953\u00f78=119, 1\",\n    \"664\u00f75=132, 4\"\n  ],\n  [\n    \"799\u00f75=159, 4\",\n    \"895\u00f72=447, 1\"\n  ],\n  [\n    \"595\u00f77=85, 0\",\n    \"651\u00f77=93, 0\"\n  ],\n  [\n    \"781\u00f73=260, 1\",\n    \"169\u00f77=24, 1\"\n  ],\n  [\n    \"144\u00f79=16, 0\",\n    \"116\u00f75=23, 1\"\n  ]\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: true\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# List of (old, new) text pairs for each three-digit division problem cell.\n# Every old value is unique across the document, so Find/Replace All is unambiguous.\n$replacements = @(\n    @(\"473\u00f75=94, 3\", \"162\u00f79=18, 0\"),\n    @(\"758\u00f72=379, 0\", \"855\u00f75=171, 0\"),\n    @(\"310\u00f78=38, 6\", \"365\u00f78=45, 5\"),\n    @(\"128\u00f77=18, 2\", \"965\u00f73=321, 2\"),\n    @(\"607\u00f76=101, 1\", \"421\u00f78=52, 5\"),\n    @(\"971\u00f79=107, 8\", \"754\u00f79=83, 7\"),\n    @(\"698\u00f75=139, 3\", \"175\u00f76=29, 1\"),\n    @(\"240\u00f73=80, 0\", \"374\u00f79=41, 5\"),\n    @(\"764\u00f72=382, 0\", \"572\u00f78=71, 4\"),\n    @(\"919\u00f79=102, 1\", \"495\u00f79=55, 0\"),\n    @(\"219\u00f74=54, 3\", \"567\u00f75=113, 2\"),\n    @(\"538\u00f77=76, 6\", \"196\u00f78=24, 4\"),\n    @(\"399\u00f75=79, 4\", \"744\u00f79=82, 6\"),\n    @(\"972\u00f75=194, 2\", \"810\u00f79=90, 0\"),\n    @(\"770\u00f75=154, 0\", \"341\u00f77=48, 5\"),\n    @(\"701\u00f76=116, 5\", \"928\u00f76=154, 4\"),\n    @(\"421\u00f74=105, 1\", \"838\u00f79=93, 1\"),\n    @(\"987\u00f78=123, 3\", \"778\u00f73=259, 1\"),\n    @(\"237\u00f73=79, 0\", \"682\u00f78=85, 2\"),\n    @(\"326\u00f76=54, 2\", \"489\u00f78=61, 1\"),\n    @(\"953\u00f78=119, 1\", \"664\u00f75=132, 4\"),\n    @(\"799\u00f75=159, 4\", \"895\u00f72=447, 1\"),\n    @(\"595\u00f77=85, 0\", \"651\u00f77=93, 0\"),\n    @(\"781\u00f73=260, 1\", \"169\u00f77=24, 1\"),\n    @(\"144\u00f79=16, 0\", \"116\u00f75=23, 1\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null  # 2 = wdReplaceAll\n}\n\n"}
